$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same rows of data that need
# their "想去人数" (F column) values updated.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 572
    $ws.Range("F7").Value = 32
    $ws.Range("F8").Value = 488
    $ws.Range("F9").Value = 3618
}
